$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.502.52"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "3.067.32"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").Value = "3.067.52"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000237"
$ws.Range("D14").Style = "Normal"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "3.575.15"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "63.457.03"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "3.068.89"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "491.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("E27").Value = "  +11.01%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("D36").Value = "0.0₃0824"
$ws.Range("E36").Value = "  -2.68%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.59%  "
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "438.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "2.842.76"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
